# test_twice.xlsx edit
# Commit message: "deleted comments in convert_to_datetime(), implemented
# BETWEEN in bef_betw_aft()"
#
# The underlying spreadsheet is used by a script (convert_to_datetime /
# bef_betw_aft) that reads each availability cell and classifies it as
# "before", "between" or "after" a given time window. This edit:
#   - Fixes a cell on the "Sheet2" tab that said "Not Available" but should
#     read the shorthand "not avai" (matching the other "not available" /
#     "not avai" style entries used elsewhere as sentinel text).
#   - Updates the "Free except 16h45-18h15" entry to use the clearer
#     "Free between 16h45-18h15" phrasing, now that BETWEEN handling has
#     been implemented in bef_betw_aft().

$wb = $excel.ActiveWorkbook

# NOTE: workbook tab names vs. internal file order: the tab named
# "Sheet2" is the first sheet (dimension A1:J9) and is the active /
# selected tab; the tab named "Sheet1" is the second sheet
# (dimension A1:I12).
$wsSheet2 = $wb.Worksheets.Item("Sheet2")
$wsSheet1 = $wb.Worksheets.Item("Sheet1")

# --- Content changes on "Sheet2" ---
# J7: "Free except 16h45-18h15" -> "Free between 16h45-18h15"
$wsSheet2.Range("J7").Value = "Free between 16h45-18h15"

# H3: "Not Available" -> "not avai"
$wsSheet2.Range("H3").Value = "not avai"

# --- Selection / cursor bookkeeping (matches the saved sheetView state) ---
# Update the non-active sheet's cached selection first...
$null = $wsSheet1.Range("D28").Select()
# ...then finish on "Sheet2" so it remains the active/visible tab.
$null = $wsSheet2.Range("J10").Select()
